$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Store date/time as text values (zero-padded time, hashed-looking date/time strings)
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "20230906"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0900"

# Clear the now-unused columns G (location) and I (goalie list) on row 2
$ws.Range("G2").Value = $null
$ws.Range("I2").Value = $null

# Delete rows 3 and 4 (extra game entries removed)
$ws.Rows("3:4").Delete()
